# Arbeitsstunden recovery.xlsx - apply "3d Modelle und Überarbeitung" commit
# Targets sheet "Roman (November)" (2nd worksheet / sheet2.xml)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Rows 20-30: only the end-time (column C) was corrected; start time, date and
# description stay the same. Column D/F/G are formulas and recalc automatically.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 0.51041666666666663
$ws.Range("C21").Value = 0.60416666666666663
$ws.Range("C22").Value = 0.75
$ws.Range("C23").Value = 0.71875
$ws.Range("C24").Value = 0.57291666666666663
$ws.Range("C25").Value = 0.6875
$ws.Range("C26").Value = 0.38541666666666669
$ws.Range("C27").Value = 0.70138888888888884
$ws.Range("C28").Value = 0.45833333333333331
$ws.Range("C29").Value = 0.52430555555555558
$ws.Range("C30").Value = 0.53125

# ---------------------------------------------------------------------------
# Row 31: date + end time corrected
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = 44187
$ws.Range("C31").Value = 0.69791666666666663

# ---------------------------------------------------------------------------
# Row 32: date, start + end time corrected, and a description is now present
# (shared string index 41 in the final file)
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = 44192
$ws.Range("B32").Value = 0.47916666666666669
$ws.Range("C32").Value = 0.55208333333333337
$ws.Range("E32").Value = "Labborhandout: Pumpenübersicht, Füllkörper, Beschriftung mit Edding"

# ---------------------------------------------------------------------------
# Row 33: date, start + end time corrected, and a description is now present
# (shared string index 42 in the final file)
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = 44192
$ws.Range("B33").Value = 0.61458333333333337
$ws.Range("C33").Value = 0.65625
$ws.Range("E33").Value = "Laborhandout: Normschliffe, Schlifffett, Eismaschine, Ultraschallbad"

# ---------------------------------------------------------------------------
# New row 34
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = 43832
$ws.Range("B34").Value = 0.61458333333333337
$ws.Range("C34").Value = 0.71875
$ws.Range("D34").Formula = "=HOUR(C34)+MINUTE(C34)/60-HOUR(B34)-MINUTE(B34)/60+D33"
$ws.Range("E34").Value = "Laborhandout: Überarbeitung - 3D-Modell Muffe mit Erklärung"
$ws.Range("F34").Formula = "=30+`$F`$19-D34"
$ws.Range("G34").Formula = "=120-(`$D`$19+D34)"

# copy the number formats / alignment / borders from row 33 so the new row
# matches the look of the rest of the table
$ws.Range("A33:C33").Copy() | Out-Null
$ws.Range("A34:C34").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("F33:G33").Copy() | Out-Null
$ws.Range("F34:G34").PasteSpecial(-4122) | Out-Null

# re-apply values/formulas (PasteSpecial of formats only shouldn't disturb
# them, but make sure nothing was clobbered)
$ws.Range("A34").Value = 43832
$ws.Range("B34").Value = 0.61458333333333337
$ws.Range("C34").Value = 0.71875
$ws.Range("D34").Formula = "=HOUR(C34)+MINUTE(C34)/60-HOUR(B34)-MINUTE(B34)/60+D33"
$ws.Range("F34").Formula = "=30+`$F`$19-D34"
$ws.Range("G34").Formula = "=120-(`$D`$19+D34)"

# ---------------------------------------------------------------------------
# New row 35 (wraps to two lines, like row 27)
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = 43832
$ws.Range("B35").Value = 0.51041666666666663
$ws.Range("C35").Value = 0.64583333333333337
$ws.Range("D35").Formula = "=HOUR(C35)+MINUTE(C35)/60-HOUR(B35)-MINUTE(B35)/60+D34"
$ws.Range("E35").Value = "Laborhandout: Überarbeitung Layout +  Formulierung, `nRecherche und Erstellung von 3D Modellen für Rührer"
$ws.Range("F35").Formula = "=30+`$F`$19-D35"
$ws.Range("G35").Formula = "=120-(`$D`$19+D35)"

$ws.Range("A33:C33").Copy() | Out-Null
$ws.Range("A35:C35").PasteSpecial(-4122) | Out-Null
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4122) | Out-Null
$ws.Range("F33:G33").Copy() | Out-Null
$ws.Range("F35:G35").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4122) | Out-Null

$ws.Range("A35").Value = 43832
$ws.Range("B35").Value = 0.51041666666666663
$ws.Range("C35").Value = 0.64583333333333337
$ws.Range("D35").Formula = "=HOUR(C35)+MINUTE(C35)/60-HOUR(B35)-MINUTE(B35)/60+D34"
$ws.Range("E35").Value = "Laborhandout: Überarbeitung Layout +  Formulierung, `nRecherche und Erstellung von 3D Modellen für Rührer"
$ws.Range("F35").Formula = "=30+`$F`$19-D35"
$ws.Range("G35").Formula = "=120-(`$D`$19+D35)"

$ws.Range("E35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 30

# ---------------------------------------------------------------------------
# Selection / view state
# ---------------------------------------------------------------------------
$ws.Range("C35").Select()
